$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 136, shifting the existing data
# (old rows 136:148) down to 138:150.
$ws.Rows("136:137").Insert()

# New row 136 - latest "Primera" quality weekly entry
$ws.Cells.Item(136, 1).Value = 9
$ws.Cells.Item(136, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(136, 3).Value = "Metropolitana"
$ws.Cells.Item(136, 4).Value = 44504
$ws.Cells.Item(136, 5).Value = 13
$ws.Cells.Item(136, 6).Value = 100112017
$ws.Cells.Item(136, 7).Value = "Apio"
$ws.Cells.Item(136, 8).Value = "Americana (o)"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 80
$ws.Cells.Item(136, 11).Value = 7000
$ws.Cells.Item(136, 12).Value = 8000
$ws.Cells.Item(136, 13).Value = 7500
$ws.Cells.Item(136, 14).Value = "$/docena de matas"
$ws.Cells.Item(136, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(136, 16).Value = 1250
$ws.Cells.Item(136, 17).Value = 6
$ws.Cells.Item(136, 18).Value = "Hortaliza"

# New row 137 - latest "Segunda" quality weekly entry
$ws.Cells.Item(137, 1).Value = 9
$ws.Cells.Item(137, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(137, 3).Value = "Metropolitana"
$ws.Cells.Item(137, 4).Value = 44504
$ws.Cells.Item(137, 5).Value = 13
$ws.Cells.Item(137, 6).Value = 100112017
$ws.Cells.Item(137, 7).Value = "Apio"
$ws.Cells.Item(137, 8).Value = "Americana (o)"
$ws.Cells.Item(137, 9).Value = "Segunda"
$ws.Cells.Item(137, 10).Value = 34
$ws.Cells.Item(137, 11).Value = 6000
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 13).Value = 6000
$ws.Cells.Item(137, 14).Value = "$/docena de matas"
$ws.Cells.Item(137, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(137, 16).Value = 1000
$ws.Cells.Item(137, 17).Value = 6
$ws.Cells.Item(137, 18).Value = "Hortaliza"

# Ensure the date column keeps the same number format used elsewhere in column D
$ws.Range("D136:D137").NumberFormat = $ws.Range("D138").NumberFormat
